$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pasDomainText = '``` SQL   -- Ass uming the following table names and columns based on the given context  -- Table: prote ins( id, name, pas_ domain_ sequence)  -- Table: mut ations( id, protein_ id, position, type)  -- Table: cell_ lines( id, name, ic50)   -- First, join prote ins and mut ations tables using protein_ id  SELECTp. name AS protein_ name,m. position,m. type  FROM prote insp  JO IN mut ationsm ONp. id=m. prote in_ id W HEREp. name='' target_ prote in'' -- replace with the actual protein name  ANDp. pas_ domain_ sequence IS NOT NULL -- assuming pas_ domain_ sequence is stored in prote ins table   -- Then, filter mut ations that are located inP AS domain  -- This assumes that the length ofP AS domain is known and stored ina variable or hardc oded value  -- Re place'' start_ index'' and'' end_ index'' with actual values  -- For example, if theP AS domain of the target protein is between positions 100 and 200:  ANDm. positionB ET WE EN start_ index AND end_ index   -- Join the result with cell_ lines table using ic50 values  -- Ass uming the select ivity is in vers ely proport ional to ic50, so lower ic50 values indicate higher select ivity for pro state cancer cell lines  -- Re place'' pro state_ cell_ line_ id'' with the actual id of the pro state cell line  SELECT protein_ name, position, type, ic50  FROM(   SE LECTp. name AS protein_ name,m. position,m. type   FROM prote insp   JO IN mut ationsm ONp. id=m. prote in_ id   WHEREp. name='' target_ prote in''   ANDp. pas_ domain_ sequence IS NOT NULL   ANDm. positionB ET WE EN start_ index AND end_ index ) mut ations  JO IN cell_ lines cl ON mut ations. id= cl. mut ation_ id -- assuming there isa mut ation_ id column in cell_ lines table W HERE cl. id= pro state_ cell_ line_ id  ORDER BY'

$ccleBlobText = '```sql--JointheCCLE_GisticCNDBandCCLEMutDBtablesontheDepMap_IDcolumnforgenePOLR3ESELECTc.CCLE_Name,g.gistic_cnFROMccle_cn_gisticcJOINccle_mutationgONc.DepMap_ID=g.DepMap_IDWHEREHugo_Symbol=''POLR3E'';--Gettheabsolutevalueofgistic_cnandcheckifitissignificantlydifferentbetweencancerandnormaltissues--Thisqueryassumesthatthereisaseparatetableorcolumnthatindicatesthetissuetype(cancerornormal)foreachentryintheCCLE_GisticCNDBtable--Replace''tissue_type_table''and''tissue_type_column''withtheactualnamesofthetableandcolumnthatcontainthisinformation--Thisqueryalsoassumesthatthereisafunctionorwaytodetermineifthegeneisoverexpressedorunderexpressedbasedonthegistic_cnvalueDECLARE@POLR3E_gistic_cnVARCHAR(10);SET@POLR3E_gistic_cn=(SELECTABS(CAST(gistic_cnASINT))FROMccle_cn_gisticJOINccle_mutationmONccle_cn_gistic.DepMap_ID=ccle_mutation.DepMap_IDWHEREHugo_Symbol=''POLR3E'');SELECT''POLR3E''ASGene_Symbol,@POLR3E_gistic_cnASgistic_cn,(--ReplacethisquerywiththeactuallogictodetermineifthegeneisoverexpressedorunderexpressedCASEWHEN(SELECTCOUNT(*)FROMtissue_type_tableWHEREtissue_type=''cancer''ANDDepMap_ID=c.DepMap_IDANDABS(CAST(gistic_cnASINT))>(--ReplacethisquerywiththeactuallogictodeterminethethresholdforsignificantdifferentialexpressionSELECTAVG(ABS(CAST(gistic_cnASINT)))FROMtissue_type_tableWHEREtissue'

$ws.Range("C2").Value = $pasDomainText
$ws.Range("C3").Value = $ccleBlobText
$ws.Range("C4").Value = $ccleBlobText
$ws.Range("C5").Value = $ccleBlobText
$ws.Range("C6").Value = $ccleBlobText
$ws.Range("C7").Value = $ccleBlobText
